$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Add new personnel row values for Zoe Sandwith (row 6)
# Shared-string insertion order follows the order in which new values are
# first written, so match the target uniqueCount ordering: ORCID, then
# middle initial, then email.
$ws.Range("F6").Value = "0000-0001-9952-9526"
$ws.Range("B6").Value = "O"
$ws.Range("E6").Value = "zoe.sandwith@hakai.org"

# Update the active selection on the sheet to match the final author selection
$ws.Range("E6").Select()
